$p = $ppt.ActivePresentation

# --- Slide 1: resize/position the "Subtitle 4" placeholder (was inheriting
#     the layout's default <p:spPr/>, now gets an explicit xfrm - shape was
#     dragged taller) ---
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$subtitle.Left = 250.20723724365234
$subtitle.Top = 164.24889373779297
$subtitle.Width = 652.9383239746094
$subtitle.Height = 261.62339782714844

# --- Handout master & notes master: bump the auto date placeholder text
#     from 4/6/2021 to 4/8/21 ---
$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.Text = "4/8/21"

$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "4/8/21"
